$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing data rows (2-10) with new forecast-error values ---

# Row 2 (Q0 -> index 6)
$ws.Range("B2").Value = 0.1578884585375505
$ws.Range("C2").Value = 1.093026323637141
$ws.Range("D2").Value = 3.780296758266214
$ws.Range("E2").Value = 1.944298526015543
$ws.Range("F2").Value = 1.957160022273444
$ws.Range("G2").Value = 51

# Row 3 (Q1)
$ws.Range("B3").Value = 0.3634719299439699
$ws.Range("C3").Value = 1.11718029823843
$ws.Range("D3").Value = 3.396285412598066
$ws.Range("E3").Value = 1.842901357261985
$ws.Range("F3").Value = 1.825045028446704
$ws.Range("G3").Value = 50

# Row 4 (Q2)
$ws.Range("B4").Value = 0.2069889977446989
$ws.Range("C4").Value = 1.057700982088954
$ws.Range("D4").Value = 3.066631775569693
$ws.Range("E4").Value = 1.751180109403283
$ws.Range("F4").Value = 1.756924272632518
$ws.Range("G4").Value = 49

# Row 5 (Q3)
$ws.Range("B5").Value = 0.3612381420177023
$ws.Range("C5").Value = 1.238224076027852
$ws.Range("D5").Value = 3.634105439593412
$ws.Range("E5").Value = 1.906332982349467
$ws.Range("F5").Value = 1.891601805341011
$ws.Range("G5").Value = 48

# Row 6 (Q4)
$ws.Range("B6").Value = 0.2411319408316268
$ws.Range("C6").Value = 1.211902555076712
$ws.Range("D6").Value = 3.480821888571079
$ws.Range("E6").Value = 1.865696086872425
$ws.Range("F6").Value = 1.870049010965773
$ws.Range("G6").Value = 47

# Row 7 (Q5)
$ws.Range("B7").Value = 0.3837287414956624
$ws.Range("C7").Value = 1.254588007131292
$ws.Range("D7").Value = 3.75232404806467
$ws.Range("E7").Value = 1.937091646790278
$ws.Range("F7").Value = 1.919684637565255
$ws.Range("G7").Value = 46

# Row 8 (Q6)
$ws.Range("B8").Value = 0.2264048185344497
$ws.Range("C8").Value = 1.210514367268617
$ws.Range("D8").Value = 3.438743821690994
$ws.Range("E8").Value = 1.854385025201345
$ws.Range("F8").Value = 1.861309476689316
$ws.Range("G8").Value = 45

# Row 9 (Q7)
$ws.Range("B9").Value = 0.3609039862410751
$ws.Range("C9").Value = 1.204350968061618
$ws.Range("D9").Value = 3.365592077467562
$ws.Range("E9").Value = 1.834555008024442
$ws.Range("F9").Value = 1.81950016883021
$ws.Range("G9").Value = 44

# Row 10 (Q8) - also gains a previously-missing F value
$ws.Range("B10").Value = 0.2857062726838135
$ws.Range("C10").Value = 1.290388063109146
$ws.Range("D10").Value = 3.826836026517005
$ws.Range("E10").Value = 1.956230054599153
$ws.Range("F10").Value = 1.958157187300658
$ws.Range("G10").Value = 43

# --- Add new row 11 (Q9) ---
# Copy formatting/style from row 10's label cell so the new label cell matches
# the existing bold/centered/bordered header-like style used in column A.
$ws.Range("A10").Copy($ws.Range("A11"))
$ws.Range("A11").Value = "Q9"
$ws.Range("B11").Value = 0.3026413976930326
$ws.Range("C11").Value = 1.355495109822078
$ws.Range("D11").Value = 3.884252099312994
$ws.Range("E11").Value = 1.970850602991763
$ws.Range("F11").Value = 1.971081985376132
$ws.Range("G11").Value = 42
